# Add the new dialog system.
# - Inserts a new "dialogTag" column (C) into the dialog table on sheet "第二章"
# - Rewrites the sample dialogue rows with a new 4-line conversation
# - Clears the now-unused Sheet2 scratch table

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("第二章")

# --- Insert a new column C ("dialogTag"). Everything right of B shifts from
#     C..H to D..I. Match the format/width of the (now-shifted) neighbouring
#     "name" column so the new column looks like part of the same table.
$ws1.Columns("C").Insert()
$ws1.Columns("C").ColumnWidth = $ws1.Columns("D").ColumnWidth
$ws1.Range("D1:D7").Copy()
$ws1.Range("C1:C7").PasteSpecial(-4122)

# --- Row 1: column headers
$ws1.Range("C1").Value = "dialogTag"

# --- Row 2: ##type row - the new column is also a plain string
$ws1.Range("C2").Value = "string"

# --- Row 3: ## comment row describing each column
$ws1.Range("C3").Value = "对每一组对话的标识符"

# --- Rows 4-7: the actual dialogue data, replaced with the new conversation
$ws1.Range("B4").Value = "D00001"
$ws1.Range("C4").Value = "MainChapter2"
$ws1.Range("D4").Value = "洛斯"
$ws1.Range("E4").Value = "喂，是小兔子吗？我是之前选拔赛评委席上的洛斯，你还记得我吗？"
$ws1.Range("F4").Value = "D00002"

$ws1.Range("B5").Value = "D00002"
$ws1.Range("C5").Value = "MainChapter2"
$ws1.Range("D5").Value = "#name#"
$ws1.Range("E5").Value = "呃，我叫#name#，记得……您是那个尾巴总甩个不停的……嗯，评委？"
$ws1.Range("F5").Value = "D00003"

$ws1.Range("B6").Value = "D00003"
$ws1.Range("C6").Value = "MainChapter2"
$ws1.Range("D6").Value = "洛斯"
$ws1.Range("E6").Value = "哈哈，那就好。我对你印象挺深，想和你聊聊。明天有空吗？"
$ws1.Range("F6").Value = "D00004"

$ws1.Range("B7").Value = "D00004"
$ws1.Range("C7").Value = "MainChapter2"
$ws1.Range("D7").Value = "#name#"
$ws1.Range("E7").Value = "哦？难道是关于比赛的事情吗?"
$ws1.Range("F7").Value = -1

# Rows 4-6 grew a touch taller to match row 7's wrapped text.
$ws1.Rows("4:6").RowHeight = 28.8

# Selection marker matches the authored edit.
$ws1.Range("E9").Select()

# --- Sheet2 was a leftover scratch copy of the table; the new dialog data
#     now lives solely on sheet 1, so clear it out entirely.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows("1:6").Delete()
